$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "68496"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "2534"
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "52565"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "3042"
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "55846"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "2880"
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "58972"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "2763"
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "65038"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "2597"
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "65680"
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "67420"
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "96803"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "1700"
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "101322"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "59231345"
$ws.Cells.Item(10, 3).Value = "Player-59231345"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "1602"
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "102250"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "44437839"
$ws.Cells.Item(11, 3).Value = "strangetamer828"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "1587"
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "113592"
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "114728"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "1456"
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "36127"
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "10919"
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = "53060417"
$ws.Cells.Item(19, 3).Value = "㊥老纳信耶稣"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "5894"
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "10716"
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = "49710892"
$ws.Cells.Item(20, 3).Value = "MMMMMMM"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "5912"
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "13882"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "5605"
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "18483"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "5287"
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "18780"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "5271"
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "23027"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "5019"
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "39859"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "4162"
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "59895"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "2733"
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "64849"
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "14835"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "5527"
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "16518"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "5406"
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "16649"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "5398"
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "17899"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "5318"
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "21477"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "5111"
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "22219"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "5066"
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "24657"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "4921"
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "28146"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "4740"
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "30960"
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "3649043"
$ws.Cells.Item(36, 3).Value = "Dj6106"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "4611"
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "33100"
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = "56732705"
$ws.Cells.Item(37, 3).Value = "时间温柔皆遗憾"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "4515"
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "33547"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "4495"
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "34026"
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = "56585361"
$ws.Cells.Item(39, 3).Value = "`"㊥ go策划我要ali`""
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "4473"
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "38374"
$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "58408326"
$ws.Cells.Item(40, 3).Value = "`"Killer Bee`""
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "4243"
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "40777"
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "1304123"
$ws.Cells.Item(41, 3).Value = "Cccccccccccc"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "4110"
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = "6598"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "6302"
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = "9035"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "6047"
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "13673"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "5625"
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "13743"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "5618"
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = "14481"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "5554"
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = "17141"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "5365"
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = "17662"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "5333"
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = "20312"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "5184"
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = "24863"
$ws.Cells.Item(51, 1).NumberFormat = "@"
$ws.Cells.Item(51, 1).Value = "30397"
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = "32316256"
$ws.Cells.Item(51, 3).Value = "`"秋の風 ..`""
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "4636"
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = "30416"
$ws.Cells.Item(52, 2).NumberFormat = "@"
$ws.Cells.Item(52, 2).Value = "47459684"
$ws.Cells.Item(52, 3).Value = "㊥阿闹切克闹"
$ws.Cells.Item(52, 5).NumberFormat = "@"
$ws.Cells.Item(52, 5).Value = "4635"
$ws.Cells.Item(53, 1).NumberFormat = "@"
$ws.Cells.Item(53, 1).Value = "33414"
$ws.Cells.Item(53, 2).NumberFormat = "@"
$ws.Cells.Item(53, 2).Value = "56573048"
$ws.Cells.Item(53, 3).Value = "Xiaotian"
$ws.Cells.Item(53, 5).NumberFormat = "@"
$ws.Cells.Item(53, 5).Value = "4502"
$ws.Cells.Item(54, 1).NumberFormat = "@"
$ws.Cells.Item(54, 1).Value = "33703"
$ws.Cells.Item(54, 2).NumberFormat = "@"
$ws.Cells.Item(54, 2).Value = "56379103"
$ws.Cells.Item(54, 3).Value = "Globalking"
$ws.Cells.Item(54, 5).NumberFormat = "@"
$ws.Cells.Item(54, 5).Value = "4488"
$ws.Cells.Item(55, 1).NumberFormat = "@"
$ws.Cells.Item(55, 1).Value = "33874"
$ws.Cells.Item(55, 2).NumberFormat = "@"
$ws.Cells.Item(55, 2).Value = "37069173"
$ws.Cells.Item(55, 3).Value = "詹toniii"
$ws.Cells.Item(55, 5).NumberFormat = "@"
$ws.Cells.Item(55, 5).Value = "4481"
$ws.Cells.Item(56, 1).NumberFormat = "@"
$ws.Cells.Item(56, 1).Value = "36565"
$ws.Cells.Item(56, 5).NumberFormat = "@"
$ws.Cells.Item(56, 5).Value = "4340"
$ws.Cells.Item(57, 1).NumberFormat = "@"
$ws.Cells.Item(57, 1).Value = "37228"
$ws.Cells.Item(57, 2).NumberFormat = "@"
$ws.Cells.Item(57, 2).Value = "38893233"
$ws.Cells.Item(57, 3).Value = "`"快乐 二哈`""
$ws.Cells.Item(57, 5).NumberFormat = "@"
$ws.Cells.Item(57, 5).Value = "4304"
$ws.Cells.Item(58, 1).NumberFormat = "@"
$ws.Cells.Item(58, 1).Value = "40264"
$ws.Cells.Item(58, 5).NumberFormat = "@"
$ws.Cells.Item(58, 5).Value = "4138"
$ws.Cells.Item(59, 1).NumberFormat = "@"
$ws.Cells.Item(59, 1).Value = "41837"
$ws.Cells.Item(59, 5).NumberFormat = "@"
$ws.Cells.Item(59, 5).Value = "4041"
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "41880"
$ws.Cells.Item(60, 2).NumberFormat = "@"
$ws.Cells.Item(60, 2).Value = "57813281"
$ws.Cells.Item(60, 3).Value = "XAUEN"
$ws.Cells.Item(60, 5).NumberFormat = "@"
$ws.Cells.Item(60, 5).Value = "4039"
$ws.Cells.Item(61, 1).NumberFormat = "@"
$ws.Cells.Item(61, 1).Value = "42749"
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = "55634661"
$ws.Cells.Item(61, 3).Value = "Opalus"
$ws.Cells.Item(61, 5).NumberFormat = "@"
$ws.Cells.Item(61, 5).Value = "3992"
$ws.Cells.Item(62, 1).NumberFormat = "@"
$ws.Cells.Item(62, 1).Value = "43420"
$ws.Cells.Item(62, 2).NumberFormat = "@"
$ws.Cells.Item(62, 2).Value = "59020292"
$ws.Cells.Item(62, 3).Value = "Sharnoth"
$ws.Cells.Item(62, 5).NumberFormat = "@"
$ws.Cells.Item(62, 5).Value = "3934"
$ws.Cells.Item(63, 1).NumberFormat = "@"
$ws.Cells.Item(63, 1).Value = "46807"
$ws.Cells.Item(63, 5).NumberFormat = "@"
$ws.Cells.Item(63, 5).Value = "3527"
$ws.Cells.Item(64, 1).NumberFormat = "@"
$ws.Cells.Item(64, 1).Value = "56723"
$ws.Cells.Item(64, 5).NumberFormat = "@"
$ws.Cells.Item(64, 5).Value = "2844"
$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "54328"
$ws.Cells.Item(66, 5).NumberFormat = "@"
$ws.Cells.Item(66, 5).Value = "2952"
$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "60888"
$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "31015"
$ws.Cells.Item(70, 5).NumberFormat = "@"
$ws.Cells.Item(70, 5).Value = "4608"
$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "35294"
$ws.Cells.Item(71, 5).NumberFormat = "@"
$ws.Cells.Item(71, 5).Value = "4406"
$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = "48519"
$ws.Cells.Item(72, 5).NumberFormat = "@"
$ws.Cells.Item(72, 5).Value = "3355"
$ws.Cells.Item(73, 1).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "51962"
$ws.Cells.Item(73, 5).NumberFormat = "@"
$ws.Cells.Item(73, 5).Value = "3081"
$ws.Cells.Item(74, 1).NumberFormat = "@"
$ws.Cells.Item(74, 1).Value = "54680"
$ws.Cells.Item(74, 5).NumberFormat = "@"
$ws.Cells.Item(74, 5).Value = "2933"
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "72193"
$ws.Cells.Item(75, 5).NumberFormat = "@"
$ws.Cells.Item(75, 5).Value = "2490"
$ws.Cells.Item(76, 1).NumberFormat = "@"
$ws.Cells.Item(76, 1).Value = "91523"
$ws.Cells.Item(76, 5).NumberFormat = "@"
$ws.Cells.Item(76, 5).Value = "1972"
$ws.Cells.Item(77, 1).NumberFormat = "@"
$ws.Cells.Item(77, 1).Value = "111655"
$ws.Cells.Item(77, 2).NumberFormat = "@"
$ws.Cells.Item(77, 2).Value = "57219176"
$ws.Cells.Item(77, 3).Value = "青莲道人"
$ws.Cells.Item(77, 5).NumberFormat = "@"
$ws.Cells.Item(77, 5).Value = "1497"
$ws.Cells.Item(78, 1).NumberFormat = "@"
$ws.Cells.Item(78, 1).Value = "111889"
$ws.Cells.Item(78, 2).NumberFormat = "@"
$ws.Cells.Item(78, 2).Value = "54941706"
$ws.Cells.Item(78, 3).Value = "AlexMenjivar20"
$ws.Cells.Item(78, 5).NumberFormat = "@"
$ws.Cells.Item(78, 5).Value = "1496"
$ws.Cells.Item(79, 1).NumberFormat = "@"
$ws.Cells.Item(79, 1).Value = "133824"
$ws.Cells.Item(90, 1).NumberFormat = "@"
$ws.Cells.Item(90, 1).Value = "45123"
$ws.Cells.Item(90, 5).NumberFormat = "@"
$ws.Cells.Item(90, 5).Value = "3720"
$ws.Cells.Item(93, 1).NumberFormat = "@"
$ws.Cells.Item(93, 1).Value = "111254"
$ws.Cells.Item(94, 1).NumberFormat = "@"
$ws.Cells.Item(94, 1).Value = "200235"
$ws.Cells.Item(95, 1).NumberFormat = "@"
$ws.Cells.Item(95, 1).Value = "192541"
$ws.Cells.Item(96, 1).NumberFormat = "@"
$ws.Cells.Item(96, 1).Value = "224374"
